$wb = $excel.ActiveWorkbook

# Rename the TRPO06_InternalParticipation sheet to TRPO06_ExternalParticipation
$ws = $wb.Worksheets.Item("TRPO06_InternalParticipation")
$ws.Name = "TRPO06_ExternalParticipation"

# Delete the Participation_Discount column (column L) entirely - shifts remaining
# columns left by one
$ws.Columns("L").Delete()

# Make this sheet the active / selected sheet, and select cell J13 as in the diff
$ws.Activate()
$ws.Range("J13").Select()
